# Apply the "aspects.xlsx" edit described in the commit:
#   "More bug fixes. Maybe we can improve early aspect conflict checking in the aspect map?"
#
# This adds a LABEL axis to the Edges conflict-matrix on the "Category usage"
# sheet (both as a new row and a new column), shifting MULT_IN/MULT_OUT/ASSOC
# down/right by one, and fills in the resulting symmetric matrix. It also
# fixes Q31 (ROLE/LABEL cell) from X to c6, and nudges the window size /
# selected cell to match the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Category usage")

$xlPasteFormats = -4122

function Set-StyledText($cellRef, $donorRef, $text) {
    # Writes $text first, THEN copies formatting (cell style) from $donorRef
    # onto $cellRef -- order matters, a later Value2 write resets the style.
    $ws.Range($cellRef).Value2 = $text
    $ws.Range($donorRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats) | Out-Null
}

function Set-StyledFormula($cellRef, $donorRef, $formula) {
    # Copies formatting (cell style) from $donorRef onto $cellRef, then writes a formula.
    $ws.Range($cellRef).Formula = $formula
    $ws.Range($donorRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats) | Out-Null
}

# --- Q31: X -> c6 (style s="6" unchanged) ---
Set-StyledText "Q31" "N30" "c6"

# --- Row 39 header: insert LABEL column, shift MULT_I- / MULT_OUT / ASSOC right ---
# Styles (s="3") are unchanged on all four cells, so plain value writes suffice.
$ws.Range("N39").Value2 = "LABEL"
$ws.Range("O39").Value2 = "MULT_I-"
$ws.Range("P39").Value2 = "MULT_OUT"
$ws.Range("Q39").Value2 = "ASSOC"

# --- Row 40 (REMARK) ---
Set-StyledText "N40" "D29" "-"
Set-StyledText "Q40" "C21" "-"

# --- Row 41 (SORT) ---
Set-StyledText "N41" "D29" "-"
Set-StyledText "O41" "C21" "-"
Set-StyledText "P41" "E23" "X"
Set-StyledText "Q41" "C21" "-"

# --- Row 42: MULT_IN -> LABEL ---
$ws.Range("K42").Value2 = "LABEL"
Set-StyledText "P42" "E23" "X"
Set-StyledText "Q42" "B3" "X"

# --- Row 43: MULT_OUT -> MULT_IN ---
$ws.Range("K43").Value2 = "MULT_IN"
Set-StyledText "Q43" "E30" "c1"

# --- Row 44: ASSOC -> MULT_OUT ---
$ws.Range("K44").Value2 = "MULT_OUT"
Set-StyledText "Q44" "E23" "X"

# --- Row 45: brand-new ASSOC row ---
$ws.Range("K45").Value2 = "ASSOC"
Set-StyledFormula "L45" "L44" "=Q40"
Set-StyledFormula "M45" "L44" "=Q41"
Set-StyledFormula "N45" "L44" "=Q42"
Set-StyledFormula "O45" "L44" "=Q43"
Set-StyledFormula "P45" "L44" "=Q44"

# Q45 becomes the diagonal (self) cell: blank, shaded style like the other
# diagonal cells (e.g. L45's sibling diagonal G45/Q34).
$ws.Range("Q45").ClearContents() | Out-Null
$ws.Range("Q34").Copy() | Out-Null
$ws.Range("Q45").PasteSpecial($xlPasteFormats) | Out-Null

# Recalculate so every formula (N34 etc.) carries the refreshed cached value.
$excel.Calculate()

# --- View state: selected cell + window size ---
$ws.Range("Q32").Select() | Out-Null
$wb.Windows.Item(1).WindowState = -4143  # xlNormal, harmless if unsupported
$wb.Windows.Item(1).Width = 38640 / 2
$wb.Windows.Item(1).Height = 21120 / 2
